$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "En propiedad"
$ws.Range("A2").Value = "En alquiler"
$ws.Range("A3").Value = "Otra forma"
$ws.Range("A4").Value = "Cedida gratis o a bajo precio por otro hogar, la empresa,¿"
